$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 11 header row ---
$ws.Range("A53").Value = "Week 11"
$ws.Range("A53").Font.Bold = $true

# --- Daily entries for Week 11 (rows 54-58) ---
$dates = @(42284, 42285, 42286, 42288, 42289)
$from  = @(0, 0.91666666666666663, 0.375, 0.625, 0.375)
$to    = @(0.16666666666666666, 0.08333333333333333, 0.91666666666666663, 0.08333333333333333, 0.08333333333333333)
$brk   = @(0, 0, 3, 3, 12)

for ($i = 0; $i -lt 5; $i++) {
    $r = 54 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 1).NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'

    $ws.Cells.Item($r, 2).Value = $from[$i]
    $ws.Cells.Item($r, 2).NumberFormat = "h:mm AM/PM"

    $ws.Cells.Item($r, 3).Value = $to[$i]
    $ws.Cells.Item($r, 3).NumberFormat = "h:mm AM/PM"

    $ws.Cells.Item($r, 4).Value = $brk[$i]
}

# Shared "duration" formula across the week's data rows
$ws.Range("E54:E58").NumberFormat = "0.00"
$ws.Range("E54:E58").Formula = "=MOD(C54-B54,1)*24-D54"

# --- Week 11 total row (59) ---
$ws.Range("D59").Value = "Total"
$ws.Range("D59").Font.Bold = $true

$ws.Range("E59").NumberFormat = "0.00"
$ws.Range("E59").Font.Bold = $true
$ws.Range("E59").Formula = "=SUM(E54:E58)"

# --- Update view state to match where the log now ends ---
$ws.Range("F64").Select()
